$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Commit message: "adicionando os testes negativos" (adding negative tests)
# Data change: D2 now holds the test result value "Pass"
$ws.Range("D2").Value = "Pass"
